$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I1").Value = "clado"
$ws.Range("J1").Value = "linaje"
$ws.Range("I2").Value = "23G"
$ws.Range("J2").Value = "GK.1.6"
$ws.Range("I3").Value = "23D"
$ws.Range("J3").Value = "EG.5"
$ws.Range("I4").Value = "23D"
$ws.Range("J4").Value = "EG.5.2"
$ws.Range("I5").Value = "23D"
$ws.Range("J5").Value = "EG.5.2"
$ws.Range("I6").Value = "23D"
$ws.Range("J6").Value = "EG.5.2"
$ws.Range("I7").Value = "23D"
$ws.Range("J7").Value = "EG.5.2"
$ws.Range("I8").Value = "23H"
$ws.Range("J8").Value = "HK.3"
$ws.Range("I9").Value = "23F"
$ws.Range("J9").Value = "HV.1"
$ws.Range("I10").Value = "23F"
$ws.Range("J10").Value = "HV.1"
$ws.Range("I11").Value = "23F"
$ws.Range("J11").Value = "HV.1"
$ws.Range("I12").Value = "23A"
$ws.Range("J12").Value = "GN.1.1"
$ws.Range("I13").Value = "23G"
$ws.Range("J13").Value = "GK.1.6.1"
$ws.Range("I14").Value = "23F"
$ws.Range("J14").Value = "EG.5.1.1"
$ws.Range("I15").Value = "23D"
$ws.Range("J15").Value = "FL.1.5.1"
$ws.Range("I16").Value = "23F"
$ws.Range("J16").Value = "HV.1"
$ws.Range("I17").Value = "23D"
$ws.Range("J17").Value = "EG.5.2"
$ws.Range("I18").Value = "23A"
$ws.Range("J18").Value = "JD.1.1"
$ws.Range("I19").Value = "23F"
$ws.Range("J19").Value = "HV.1"
$ws.Range("I20").Value = "recombinant"
$ws.Range("J20").Value = "XCL"
$ws.Range("I21").Value = "23F"
$ws.Range("J21").Value = "HK.13"
$ws.Range("I22").Value = "23A"
$ws.Range("J22").Value = "XBB.1.5.72"
$ws.Range("I23").Value = "23F"
$ws.Range("J23").Value = "HV.1.1"
$ws.Range("I24").Value = "23D"
$ws.Range("J24").Value = "FL.15.1.1"
$ws.Range("I25").Value = "23F"
$ws.Range("J25").Value = "HK.26"
$ws.Range("I26").Value = "23A"
$ws.Range("J26").Value = "JD.1.1.1"
$ws.Range("I27").Value = "23F"
$ws.Range("J27").Value = "HV.1"
$ws.Range("I28").Value = "23F"
$ws.Range("J28").Value = "HK.26"
$ws.Range("I29").Value = "23A"
$ws.Range("J29").Value = "JD.1.1"
$ws.Range("I30").Value = "23D"
$ws.Range("J30").Value = "EG.5.2"
$ws.Range("I31").Value = "23D"
$ws.Range("J31").Value = "EG.5.2"
$ws.Range("I32").Value = "23A"
$ws.Range("J32").Value = "GN.1"
$ws.Range("I33").Value = "23F"
$ws.Range("J33").Value = "JG.3"
$ws.Range("I34").Value = "23F"
$ws.Range("J34").Value = "EG.5.1"
$ws.Range("I35").Value = "23F"
$ws.Range("J35").Value = "HV.1"
$ws.Range("I36").Value = "23D"
$ws.Range("J36").Value = "KC.1"
$ws.Range("I37").Value = "23A"
$ws.Range("J37").Value = "FD.5.1"
$ws.Range("I38").Value = "23F"
$ws.Range("J38").Value = "HV.1"
$ws.Range("I39").Value = "23F"
$ws.Range("J39").Value = "JG.3"
$ws.Range("I40").Value = "23G"
$ws.Range("J40").Value = "GK.1.6.1"
$ws.Range("I41").Value = "23D"
$ws.Range("J41").Value = "FL.15.1.1"
$ws.Range("I42").Value = "23I"
$ws.Range("J42").Value = "JN.1"
$ws.Range("I43").Value = "23A"
$ws.Range("J43").Value = "XBB.1.5.72"
$ws.Range("I44").Value = "23I"
$ws.Range("J44").Value = "JN.1.1"
$ws.Range("I45").Value = "23H"
$ws.Range("J45").Value = "HK.3.2"
$ws.Range("I46").Value = "23I"
$ws.Range("J46").Value = "JN.1"
$ws.Range("I47").Value = "23F"
$ws.Range("J47").Value = "HK.26"
$ws.Range("I48").Value = "23A"
$ws.Range("J48").Value = "XBB.1.5.109"
$ws.Range("I49").Value = "23G"
$ws.Range("J49").Value = "GK.1.6.1"
$ws.Range("I50").Value = "23D"
$ws.Range("J50").Value = "EG.5"
$ws.Range("I51").Value = "23F"
$ws.Range("J51").Value = "JG.3"
$ws.Range("I52").Value = "23A"
$ws.Range("J52").Value = "JD.1.2"
$ws.Range("I53").Value = "23A"
$ws.Range("J53").Value = "FD.5.1"
$ws.Range("I54").Value = "23F"
$ws.Range("J54").Value = "EG.5.1.3"
$ws.Range("I55").Value = "23F"
$ws.Range("J55").Value = "EG.5.1.3"
$ws.Range("I56").Value = "23F"
$ws.Range("J56").Value = "EG.5.1.3"
$ws.Range("I57").Value = "23F"
$ws.Range("J57").Value = "EG.5.1.3"
$ws.Range("I58").Value = "23F"
$ws.Range("J58").Value = "EG.5.1"
$ws.Range("I59").Value = "23E"
$ws.Range("J59").Value = "XBB.2.3"
$ws.Range("I60").Value = "23A"
$ws.Range("J60").Value = "JD.1.1.1"
$ws.Range("I61").Value = "23A"
$ws.Range("J61").Value = "JD.1.1.1"
$ws.Range("I62").Value = "23A"
$ws.Range("J62").Value = "JD.1"
$ws.Range("I63").Value = "23A"
$ws.Range("J63").Value = "XBB.1.5.102"
$ws.Range("I64").Value = "23F"
$ws.Range("J64").Value = "HV.1"
$ws.Range("I65").Value = "23A"
$ws.Range("J65").Value = "XBB.1.5"
$ws.Range("I66").Value = "23I"
$ws.Range("J66").Value = "JN.1"
$ws.Range("I67").Value = "23F"
$ws.Range("J67").Value = "HV.1"
$ws.Range("I68").Value = "23A"
$ws.Range("J68").Value = "JD.1.1.1"
$ws.Range("I69").Value = "23F"
$ws.Range("J69").Value = "HV.1"
$ws.Range("I70").Value = "23A"
$ws.Range("J70").Value = "JD.1.1.1"
$ws.Range("I71").Value = "23D"
$ws.Range("J71").Value = "FL.1.5.1"
$ws.Range("I72").Value = "23B"
$ws.Range("J72").Value = "JF.1.1"
$ws.Range("I73").Value = "23D"
$ws.Range("J73").Value = "FL.1.5.1"
$ws.Range("I74").Value = "23A"
$ws.Range("J74").Value = "JD.1.1.1"
$ws.Range("I75").Value = "23A"
$ws.Range("J75").Value = "XBB.1.5"
$ws.Range("I76").Value = "23A"
$ws.Range("J76").Value = "XBB.1.5.72"
$ws.Range("I77").Value = "23F"
$ws.Range("J77").Value = "HV.1"
$ws.Range("I78").Value = "23F"
$ws.Range("J78").Value = "HV.1.1"
$ws.Range("I79").Value = "23A"
$ws.Range("J79").Value = "XBB.1.5"
$ws.Range("I80").Value = "23A"
$ws.Range("J80").Value = "JD.1.1.1"
$ws.Range("I81").Value = "23F"
$ws.Range("J81").Value = "EG.5.1.6"
$ws.Range("I82").Value = "23A"
$ws.Range("J82").Value = "JD.1.1.1"
$ws.Range("I83").Value = "23F"
$ws.Range("J83").Value = "EG.5.1.6"
$ws.Range("I84").Value = "23A"
$ws.Range("J84").Value = "XBB.1.5"
$ws.Range("I85").Value = "23D"
$ws.Range("J85").Value = "XBB.1.9.1"
$ws.Range("I86").Value = "23B"
$ws.Range("J86").Value = "JF.1.1"
$ws.Range("I87").Value = "23D"
$ws.Range("J87").Value = "FL.1.5.1"
$ws.Range("I88").Value = "23A"
$ws.Range("J88").Value = "JD.1"
$ws.Range("I89").Value = "23D"
$ws.Range("J89").Value = "EG.11"
$ws.Range("I90").Value = "23A"
$ws.Range("J90").Value = "XBB.1.5.72"
$ws.Range("I91").Value = "23A"
$ws.Range("J91").Value = "XBB.1.5.72"
$ws.Range("I92").Value = "23A"
$ws.Range("J92").Value = "XBB.1.5.72"
$ws.Range("I93").Value = "23A"
$ws.Range("J93").Value = "JD.1.1.1"
$ws.Range("I94").Value = "recombinant"
$ws.Range("J94").Value = "XDK"
$ws.Range("I95").Value = "23I"
$ws.Range("J95").Value = "JN.1"
$ws.Range("I96").Value = "23A"
$ws.Range("J96").Value = "JD.1.1.1"
$ws.Range("I97").Value = "23F"
$ws.Range("J97").Value = "HV.1"
$ws.Range("I98").Value = "23A"
$ws.Range("J98").Value = "XBB.1.5"
$ws.Range("I99").Value = "23F"
$ws.Range("J99").Value = "HV.1"
$ws.Range("I100").Value = "23F"
$ws.Range("J100").Value = "HV.1"
$ws.Range("I101").Value = "23I"
$ws.Range("J101").Value = "JN.1"
$ws.Range("I102").Value = "23I"
$ws.Range("J102").Value = "JN.1"
$ws.Range("I103").Value = "23I"
$ws.Range("J103").Value = "JN.1"
$ws.Range("I104").Value = "23A"
$ws.Range("J104").Value = "JD.1.1.1"
$ws.Range("I105").Value = "23F"
$ws.Range("J105").Value = "JG.3"
$ws.Range("I106").Value = "23E"
$ws.Range("J106").Value = "GJ.1"
$ws.Range("I107").Value = "23A"
$ws.Range("J107").Value = "HR.1.1"
